$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear out the old data rows (2-7), keeping the header row 1 intact
$ws.Range("A2:D7").ClearContents()

# Add the new single data row at row 5
$ws.Range("B5").Value = "dxgjndg,k"

# Update the selection to match the target state
$ws.Range("B5").Select()
